$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet so it lands in the
# correct final tab order (Sheet1 (2), Sheet1, Sheet2, Sheet1 (3)).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet1 (3)"

$data = New-Object 'object[,]' 29,15
$data[0,2] = 'X, Y, Z'
$data[1,0] = 1
$data[1,1] = 'Piece'
$data[1,2] = '-1.6,-1.6,-1.6'
$data[1,3] = 'Red'
$data[1,4] = 'Yellow'
$data[1,5] = 'DarkGreen'
$data[1,10] = 'D'
$data[1,11] = 'L'
$data[1,12] = 'B'
$data[1,14] = 'Bottom Back'
$data[2,0] = 2
$data[2,1] = 'Piece'
$data[2,2] = '-0.5,-1.6,-1.6'
$data[2,4] = 'Yellow'
$data[2,5] = 'DarkGreen'
$data[2,10] = 'D'
$data[2,11] = 'M'
$data[2,12] = 'B'
$data[3,0] = 3
$data[3,1] = 'Piece'
$data[3,2] = '0.6,-1.6,-1.6'
$data[3,4] = 'Yellow'
$data[3,5] = 'DarkGreen'
$data[3,6] = 'DarkOrange'
$data[3,10] = 'D'
$data[3,11] = 'R'
$data[3,12] = 'B'
$data[4,0] = 4
$data[4,1] = 'Piece'
$data[4,2] = '-1.6,-1.6,-0.5'
$data[4,3] = 'Red'
$data[4,4] = 'Yellow'
$data[4,10] = 'D'
$data[4,11] = 'L'
$data[4,12] = 'S'
$data[4,14] = 'Bottom Middle'
$data[5,0] = 5
$data[5,1] = 'Piece'
$data[5,2] = '-0.5,-1.6,-0.5'
$data[5,4] = 'Yellow'
$data[5,10] = 'D'
$data[5,11] = 'M'
$data[5,12] = 'S'
$data[6,0] = 6
$data[6,1] = 'Piece'
$data[6,2] = '0.6,-1.6,-0.5'
$data[6,4] = 'Yellow'
$data[6,6] = 'DarkOrange'
$data[6,10] = 'D'
$data[6,11] = 'R'
$data[6,12] = 'S'
$data[7,0] = 7
$data[7,1] = 'Piece'
$data[7,2] = '-1.6,-1.6,0.6'
$data[7,3] = 'Red'
$data[7,4] = 'Yellow'
$data[7,8] = 'Blue'
$data[7,10] = 'D'
$data[7,11] = 'L'
$data[7,12] = 'F'
$data[7,14] = 'Bottom Front'
$data[8,0] = 8
$data[8,1] = 'Piece'
$data[8,2] = '-0.5,-1.6,0.6'
$data[8,4] = 'Yellow'
$data[8,8] = 'Blue'
$data[8,10] = 'D'
$data[8,11] = 'M'
$data[8,12] = 'F'
$data[9,0] = 9
$data[9,1] = 'Piece'
$data[9,2] = '0.6,-1.6,0.6'
$data[9,4] = 'Yellow'
$data[9,6] = 'DarkOrange'
$data[9,8] = 'Blue'
$data[9,10] = 'D'
$data[9,11] = 'R'
$data[9,12] = 'F'
$data[10,0] = 10
$data[10,1] = 'Piece'
$data[10,2] = '-1.6,-0.5,-1.6'
$data[10,3] = 'Red'
$data[10,5] = 'DarkGreen'
$data[10,10] = 'E'
$data[10,11] = 'L'
$data[10,12] = 'B'
$data[10,14] = 'Back Middle'
$data[11,0] = 11
$data[11,1] = 'Piece'
$data[11,2] = '-0.5,-0.5,-1.6'
$data[11,5] = 'DarkGreen'
$data[11,10] = 'E'
$data[11,11] = 'M'
$data[11,12] = 'B'
$data[12,0] = 12
$data[12,1] = 'Piece'
$data[12,2] = '0.6,-0.5,-1.6'
$data[12,5] = 'DarkGreen'
$data[12,6] = 'DarkOrange'
$data[12,10] = 'E'
$data[12,11] = 'R'
$data[12,12] = 'B'
$data[13,0] = 13
$data[13,1] = 'Piece'
$data[13,2] = '-1.6,-0.5,-0.5'
$data[13,3] = 'Red'
$data[13,10] = 'E'
$data[13,11] = 'L'
$data[13,12] = 'S'
$data[13,14] = 'Middle'
$data[14,0] = 14
$data[14,1] = 'Piece'
$data[14,2] = '0.6,-0.5,-0.5'
$data[14,6] = 'DarkOrange'
$data[14,10] = 'E'
$data[14,11] = 'R'
$data[14,12] = 'S'
$data[14,14] = 'Mddle'
$data[15,0] = 15
$data[15,1] = 'Piece'
$data[15,2] = '-1.6,-0.5,0.6'
$data[15,3] = 'Red'
$data[15,8] = 'Blue'
$data[15,10] = 'E'
$data[15,11] = 'L'
$data[15,12] = 'F'
$data[15,14] = 'Front Middle'
$data[16,0] = 16
$data[16,1] = 'Piece'
$data[16,2] = '-0.5,-0.5,0.6'
$data[16,8] = 'Blue'
$data[16,10] = 'E'
$data[16,11] = 'M'
$data[16,12] = 'F'
$data[17,0] = 17
$data[17,1] = 'Piece'
$data[17,2] = '0.6,-0.5,0.6'
$data[17,6] = 'DarkOrange'
$data[17,8] = 'Blue'
$data[17,10] = 'E'
$data[17,11] = 'R'
$data[17,12] = 'F'
$data[18,0] = 18
$data[18,1] = 'Piece'
$data[18,2] = '-1.6,0.6,-1.6'
$data[18,3] = 'Red'
$data[18,5] = 'DarkGreen'
$data[18,7] = 'White'
$data[18,10] = 'U'
$data[18,11] = 'L'
$data[18,12] = 'B'
$data[18,14] = 'Top Back'
$data[19,0] = 19
$data[19,1] = 'Piece'
$data[19,2] = '-0.5,0.6,-1.6'
$data[19,5] = 'DarkGreen'
$data[19,7] = 'White'
$data[19,10] = 'U'
$data[19,11] = 'M'
$data[19,12] = 'B'
$data[20,0] = 20
$data[20,1] = 'Piece'
$data[20,2] = '0.6,0.6,-1.6'
$data[20,5] = 'DarkGreen'
$data[20,6] = 'DarkOrange'
$data[20,7] = 'White'
$data[20,10] = 'U'
$data[20,11] = 'R'
$data[20,12] = 'B'
$data[21,0] = 21
$data[21,1] = 'Piece'
$data[21,2] = '-1.6,0.6,-0.5'
$data[21,3] = 'Red'
$data[21,7] = 'White'
$data[21,10] = 'U'
$data[21,11] = 'L'
$data[21,12] = 'S'
$data[21,14] = 'Top Middle'
$data[22,0] = 22
$data[22,1] = 'Piece'
$data[22,2] = '-0.5,0.6,-0.5'
$data[22,7] = 'White'
$data[22,10] = 'U'
$data[22,11] = 'M'
$data[22,12] = 'S'
$data[23,0] = 23
$data[23,1] = 'Piece'
$data[23,2] = '0.6,0.6,-0.5'
$data[23,6] = 'DarkOrange'
$data[23,7] = 'White'
$data[23,10] = 'U'
$data[23,11] = 'R'
$data[23,12] = 'S'
$data[24,0] = 24
$data[24,1] = 'Piece'
$data[24,2] = '-1.6,0.6,0.6'
$data[24,3] = 'Red'
$data[24,7] = 'White'
$data[24,8] = 'Blue'
$data[24,10] = 'U'
$data[24,11] = 'L'
$data[24,12] = 'F'
$data[24,14] = 'Top Front'
$data[25,0] = 25
$data[25,1] = 'Piece'
$data[25,2] = '-0.5,0.6,0.6'
$data[25,7] = 'White'
$data[25,8] = 'Blue'
$data[25,10] = 'U'
$data[25,11] = 'M'
$data[25,12] = 'F'
$data[26,0] = 26
$data[26,1] = 'Piece'
$data[26,2] = '0.6,0.6,0.6'
$data[26,6] = 'DarkOrange'
$data[26,7] = 'White'
$data[26,8] = 'Blue'
$data[26,10] = 'U'
$data[26,11] = 'R'
$data[26,12] = 'F'
$data[28,2] = 'x,'

$ws.Range("A1:O29").Value = $data

# Column widths (approximate to the nearest width this engine can represent;
# input value chosen so the stored OOXML width lands as close as possible to
# the authored widths: 3.85546875, 7.5703125, 17, 12.7109375, 13.42578125,
# 13, 15, 5.85546875).
$ws.Columns.Item(1).ColumnWidth = 3.0221354166666665
$ws.Columns.Item(2).ColumnWidth = 6.736979166666667
$ws.Columns.Item(3).ColumnWidth = 16.166666666666668
$ws.Columns.Item(4).ColumnWidth = 11.877604166666666
$ws.Columns.Item(5).ColumnWidth = 11.877604166666666
$ws.Columns.Item(6).ColumnWidth = 12.592447916666666
$ws.Columns.Item(7).ColumnWidth = 12.166666666666666
$ws.Columns.Item(8).ColumnWidth = 14.166666666666666
$ws.Columns.Item(9).ColumnWidth = 14.166666666666666
$ws.Columns.Item(10).ColumnWidth = 5.022135416666667
$ws.Columns.Item(11).ColumnWidth = 5.022135416666667
$ws.Columns.Item(12).ColumnWidth = 5.022135416666667
$ws.Columns.Item(13).ColumnWidth = 5.022135416666667
$ws.Columns.Item(14).ColumnWidth = 5.022135416666667

# Make the new sheet the active tab, with C1 selected (matches the authored
# sheetView: tabSelected="1" / activeCell="C1").
$ws.Activate()
$ws.Range("C1").Select()
